# ozon fixes part 1 05.02.2026
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rows 2-7: replace with fresh order data.
#    Column A loses its banded border/fill style (becomes plain "Normal"),
#    columns B-E keep their existing (already-alternating) styles - only the
#    values change.
# ---------------------------------------------------------------------------
$newA = @(2500632, 2499894, 2498806, 2498791, 2498590, 2498138)
$newB = @(131151, 294838, 39399, 206223, 204005, 363541)
$newC = @(11.6, 11.6, 11.6, 11.6, 11.6, 11.6)
$newD = @(0, 0, 0, 0, 0, 0)
$newE = @(12, 12, 12, 12, 12, 12)

for ($i = 0; $i -lt 6; $i++) {
    $r = 2 + $i

    $ws.Range("A$r").Value = $newA[$i]
    $ws.Range("A$r").Style = "Normal"

    $ws.Range("B$r").Value = $newB[$i]
    $ws.Range("C$r").Value = $newC[$i]
    $ws.Range("D$r").Value = $newD[$i]
    $ws.Range("E$r").Value = $newE[$i]
}

# ---------------------------------------------------------------------------
# 2. Rows 8-16: the old sample rows are gone - turn them into the same blank
#    template rows used further down the sheet (copy formats from row 17).
# ---------------------------------------------------------------------------
$ws.Range("A8:E16").ClearContents()
$ws.Range("A17:E17").Copy()
for ($r = 8; $r -le 16; $r++) {
    $ws.Range("A" + $r + ":E" + $r).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Rows 18-26: column D switches from the "11" style to the plain "5"
#    style already used by the blank rows beneath (27 onward).
# ---------------------------------------------------------------------------
$ws.Range("D27").Copy()
for ($r = 18; $r -le 26; $r++) {
    $ws.Range("D" + $r).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4. Drop the trailing blank rows 94-102 (sheet now ends at row 93).
# ---------------------------------------------------------------------------
$ws.Range("A94:E102").EntireRow.Delete()

# ---------------------------------------------------------------------------
# 5. Remove the duplicate-value conditional formatting and the hyperlink.
# ---------------------------------------------------------------------------
$ws.Cells.FormatConditions.Delete()
$ws.Hyperlinks.Delete()

# ---------------------------------------------------------------------------
# 6. Update the remembered selection.
# ---------------------------------------------------------------------------
$ws.Range("C10").Select() | Out-Null
